# Commit: "took latest code changes and append my changes"
#
# Semantic change: on the "Test Cases" sheet, the Results column (D) value
# for test case row 5 (ProfileUpdateTest) flips from "PASS" to "SKIP".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")
$ws.Range("D5").Value = "SKIP"
